$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Database")

$ws.Range("B4").Value = "eychrqfbwkj65rqogd77"
$ws.Range("B5").Value = "pscale_pw_XMsJE9uSyMrK7NnecpXf71tNHNhUrMmf9sdRdIOuPim"

$ws.Range("B5").Select()
